$wb = $excel.ActiveWorkbook

# Sheet "Логин": update selection only
$wsLogin = $wb.Worksheets.Item("Логин")
$wsLogin.Select() | Out-Null
$wsLogin.Range("B15").Select() | Out-Null

# Sheet "Пароль": update cell text + selection
$wsPass = $wb.Worksheets.Item("Пароль")
$wsPass.Select() | Out-Null
$wsPass.Range("A9").Value = "1 - 7 символов, включая минимум `nодну латинскую букву и одну цифру."
$wsPass.Range("B9").Value = "1, 6, 7, 8"
$wsPass.Range("D10").Select() | Out-Null
